$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-17 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-18 Wednesday", 2) | Out-Null

$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "83-59="
$t.Cell(1,2).Range.Text = "5+38="
$t.Cell(1,3).Range.Text = "89-9="
$t.Cell(1,4).Range.Text = "7+23="
$t.Cell(1,5).Range.Text = "18+5="
$t.Cell(2,1).Range.Text = "81-21="
$t.Cell(2,2).Range.Text = "61+35="
$t.Cell(2,3).Range.Text = "27-11="
$t.Cell(2,4).Range.Text = "2+64="
$t.Cell(2,5).Range.Text = "2+90="
$t.Cell(3,1).Range.Text = "38-4="
$t.Cell(3,2).Range.Text = "92-63="
$t.Cell(3,3).Range.Text = "56-9="
$t.Cell(3,4).Range.Text = "54-34="
$t.Cell(3,5).Range.Text = "49-34="
$t.Cell(4,1).Range.Text = "40-11="
$t.Cell(4,2).Range.Text = "81-11="
$t.Cell(4,3).Range.Text = "15+15="
$t.Cell(4,4).Range.Text = "46-45="
$t.Cell(4,5).Range.Text = "58-11="
$t.Cell(5,1).Range.Text = "24-2="
$t.Cell(5,2).Range.Text = "4+29="
$t.Cell(5,3).Range.Text = "71-63="
$t.Cell(5,4).Range.Text = "13+76="
$t.Cell(5,5).Range.Text = "71+22="
$t.Cell(6,1).Range.Text = "93-13="
$t.Cell(6,2).Range.Text = "69-55="
$t.Cell(6,3).Range.Text = "82-57="
$t.Cell(6,4).Range.Text = "10+16="
$t.Cell(6,5).Range.Text = "63-54="
$t.Cell(7,1).Range.Text = "28-18="
$t.Cell(7,2).Range.Text = "80-12="
$t.Cell(7,3).Range.Text = "4+48="
$t.Cell(7,4).Range.Text = "33+21="
$t.Cell(7,5).Range.Text = "76-52="
$t.Cell(8,1).Range.Text = "71+6="
$t.Cell(8,2).Range.Text = "48-25="
$t.Cell(8,3).Range.Text = "86-3="
$t.Cell(8,4).Range.Text = "70-45="
$t.Cell(8,5).Range.Text = "92-83="
$t.Cell(9,1).Range.Text = "68-4="
$t.Cell(9,2).Range.Text = "74-70="
$t.Cell(9,3).Range.Text = "64+27="
$t.Cell(9,4).Range.Text = "35-25="
$t.Cell(9,5).Range.Text = "88-84="
$t.Cell(10,1).Range.Text = "43+4="
$t.Cell(10,2).Range.Text = "72-64="
$t.Cell(10,3).Range.Text = "94-71="
$t.Cell(10,4).Range.Text = "23-10="
$t.Cell(10,5).Range.Text = "64-41="
$t.Cell(11,1).Range.Text = "30+18="
$t.Cell(11,2).Range.Text = "11+6="
$t.Cell(11,3).Range.Text = "90+6="
$t.Cell(11,4).Range.Text = "34+28="
$t.Cell(11,5).Range.Text = "65-4="
$t.Cell(12,1).Range.Text = "52-4="
$t.Cell(12,2).Range.Text = "28+1="
$t.Cell(12,3).Range.Text = "69+9="
$t.Cell(12,4).Range.Text = "6+91="
$t.Cell(12,5).Range.Text = "10+57="
$t.Cell(13,1).Range.Text = "27+52="
$t.Cell(13,2).Range.Text = "77-49="
$t.Cell(13,3).Range.Text = "66-44="
$t.Cell(13,4).Range.Text = "36+37="
$t.Cell(13,5).Range.Text = "91-47="
$t.Cell(14,1).Range.Text = "78+19="
$t.Cell(14,2).Range.Text = "74-48="
$t.Cell(14,3).Range.Text = "92-59="
$t.Cell(14,4).Range.Text = "51+7="
$t.Cell(14,5).Range.Text = "39-26="
$t.Cell(15,1).Range.Text = "24-1="
$t.Cell(15,2).Range.Text = "43-14="
$t.Cell(15,3).Range.Text = "42+51="
$t.Cell(15,4).Range.Text = "35+14="
$t.Cell(15,5).Range.Text = "8+42="
$t.Cell(16,1).Range.Text = "10+48="
$t.Cell(16,2).Range.Text = "68+21="
$t.Cell(16,3).Range.Text = "62+27="
$t.Cell(16,4).Range.Text = "84-37="
$t.Cell(16,5).Range.Text = "24+73="
$t.Cell(17,1).Range.Text = "80-60="
$t.Cell(17,2).Range.Text = "40-34="
$t.Cell(17,3).Range.Text = "78-63="
$t.Cell(17,4).Range.Text = "27-13="
$t.Cell(17,5).Range.Text = "13+65="
$t.Cell(18,1).Range.Text = "14+69="
$t.Cell(18,2).Range.Text = "12+66="
$t.Cell(18,3).Range.Text = "10+76="
$t.Cell(18,4).Range.Text = "68-62="
$t.Cell(18,5).Range.Text = "95-36="
$t.Cell(19,1).Range.Text = "18+38="
$t.Cell(19,2).Range.Text = "46+37="
$t.Cell(19,3).Range.Text = "39+14="
$t.Cell(19,4).Range.Text = "60-54="
$t.Cell(19,5).Range.Text = "95-25="
$t.Cell(20,1).Range.Text = "91-30="
$t.Cell(20,2).Range.Text = "48-17="
$t.Cell(20,3).Range.Text = "18+69="
$t.Cell(20,4).Range.Text = "9+43="
$t.Cell(20,5).Range.Text = "88-62="
